$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Title paragraph: " [Lista de Restrições" -> " Lista de Restrições"
#    (drop the stray leading "[" that preceded the title text; the two
#    runs collapse into a single run once the bracket is gone)
# ------------------------------------------------------------------
$d.Content.Find.Execute(
    "[Lista de Restrições", $true, $false, $false, $false, $false,
    $true, 1, $false, "Lista de Restrições", 2) | Out-Null

# ------------------------------------------------------------------
# 2. Remove the "Contato com o cliente" / "Por WhatsApp e telefone."
#    row from the restrictions table entirely.
# ------------------------------------------------------------------
$table = $d.Tables.Item(1)
for ($i = $table.Rows.Count; $i -ge 1; $i--) {
    $rowText = $table.Cell($i, 1).Range.Text
    if ($rowText -like "Contato com o cliente*") {
        $table.Rows.Item($i).Delete()
    }
}

# ------------------------------------------------------------------
# 3. Spelling/grammar fix flagged by Profº Takai: "por tanto" (two
#    words) should be "portanto" (one word) in the "Espaço local"
#    justification cell.
# ------------------------------------------------------------------
$d.Content.Find.Execute(
    "por tanto na implementação", $true, $false, $false, $false, $false,
    $true, 1, $false, "portanto na implementação", 2) | Out-Null
